$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet (clears cell data, formats, and shared strings).
$ws.Cells.Clear()

# --- Write text (shared-string) cells in "first occurrence" order so the
# rebuilt sharedStrings table lands in the same index order as the target file. ---
$ws.Range("A3").Value  = "Number of articles:"
$ws.Range("A2").Value  = "GoogleAlerts Keyword:"
$ws.Range("A6").Value  = "BoardReaderKeyword:"
$ws.Range("A1").Value  = "For Google Alerts"
$ws.Range("A5").Value  = "For BoardReader"
$ws.Range("A13").Value = "For Topic Modelling"
$ws.Range("A17").Value = "Number of Questions:"
$ws.Range("A18").Value = "Question 1:"
$ws.Range("A19").Value = "Question 2:"
$ws.Range("A20").Value = "Question 3:"
$ws.Range("A21").Value = "Question 4:"
$ws.Range("A14").Value = "Name of Excel File"
$ws.Range("A9").Value  = "For Twitter"
$ws.Range("A10").Value = "Twitter Keyword:"
$ws.Range("A15").Value = "Name of Sheet:"
$ws.Range("B10").Value = "#istandwithraeesah"
$ws.Range("C3").Value  = "max = 75"
$ws.Range("C7").Value  = "max = 1000"
$ws.Range("C11").Value = "no max"
$ws.Range("B15").Value = "14-07-20 1454"
$ws.Range("B18").Value = "Text"
$ws.Range("A16").Value = "Number of Topics:"

# repeated strings (re-use existing shared-string slots)
$ws.Range("A7").Value  = "Number of articles:"
$ws.Range("A11").Value = "Number of articles:"
$ws.Range("B14").Value = "#istandwithraeesah"

# --- numeric cells ---
$ws.Range("B11").Value = 50
$ws.Range("B16").Value = 3
$ws.Range("B17").Value = 1

# --- bold section headers ---
foreach ($addr in @("A1","A5","A9","A13")) {
    $ws.Range($addr).Font.Bold = $true
}

# --- "touched" (explicitly re-applied, still-regular) font for the
# Topic-Modelling sub-header rows ---
foreach ($addr in @("A14","A15","A16")) {
    $ws.Range($addr).Font.Name = "Calibri"
    $ws.Range($addr).Font.Size = 11
}

# --- column widths (best approximation; this runtime has no real font
# metrics so AutoFit/ColumnWidth can't reproduce Excel's exact bestFit
# pixel widths) ---
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null

# --- selection / page orientation ---
$ws.Range("A6").Select() | Out-Null
$ws.PageSetup.Orientation = 1

Write-Output "done"
